# "Generate Report for Archive"
# The localization-status report is regenerated: every row whose status was
# "Ready for handoff" is now reported as "In Translation" (the handoff step
# ran and the items went back into translation), on the Overview sheet
# (zh-cn/de-de status columns E & F) as well as on each per-locale detail
# sheet (Status column C). This also removes the now-unused
# "Ready for handoff" shared string and narrows the affected status columns
# to match their new (shorter) content.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) for rows 5-7 ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($row in 5..7) {
    foreach ($col in @("E", "F")) {
        $cell = $overview.Range("$col$row")
        if ($cell.Value() -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-locale detail sheets: Status column C for rows 5-7 ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $detail = $wb.Worksheets.Item($sheetName)
    foreach ($row in 5..7) {
        $cell = $detail.Range("C$row")
        if ($cell.Value() -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
    $detail.Columns.Item(3).ColumnWidth = 12.5
}
